# Fill in the test plan for the SavingsAccount class (A02 assignment).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Developer name (was a placeholder "Student Name")
$ws.Range("C3").Value = "Raven Manalastas"

# Preconditions column for the six filled-in test cases (rows 7-12).
$ws.Range("E7").Value = "None"
$ws.Range("E8").Value = "None"
$ws.Range("E9").Value = "All inputs must be vaild"
$ws.Range("E10").Value = "All inputs must be vaild"
$ws.Range("E11").Value = "All inputs must be vaild"
$ws.Range("E12").Value = "All inputs must be vaild"

# Test case 1 (row 7): __init__ / Attributes are set to parameter values.
$ws.Range("G7").Value = "Attributes are set to input values."

# Test case 2 (row 8): __init__ / minimum_balance has invalid type.
$ws.Range("F8").Value = "minimum_balance = ""Invalid balance"""
$ws.Range("G8").Value = "set the minimum balance attribute to 50"

# Method Inputs reused across the three "Account Number = 910 ... Balance = 575 ..." cases.
$ws.Range("F7").Value = "Account Number = 910`nClient Number = 1910`nBalance = 575`nDate Created = July, 14, 2000`nminimum_balance = 100"
$ws.Range("F9").Value = "Account Number = 910`nClient Number = 1910`nBalance = 575`nDate Created = July, 14, 2000`nminimum_balance = 100"
$ws.Range("F12").Value = "Account Number = 910`nClient Number = 1910`nBalance = 575`nDate Created = July, 14, 2000`nminimum_balance = 100"

# Expected Result reused between the "balance greater than" and "balance equal to" cases.
$ws.Range("G9").Value = "setting the service_charge to 0.50"
$ws.Range("G10").Value = "setting the service_charge to 0.50"

# Test case 4 (row 10): get_service_charges / balance equal to minimum balance.
$ws.Range("F10").Value = "Account Number = 910`nClient Number = 1910`nBalance = 100`nDate Created = July, 14, 2000`nminimum_balance = 100"

# Test case 5 (row 11): get_service_charges / balance less than minimum balance.
$ws.Range("F11").Value = "Account Number = 910`nClient Number = 1910`nBalance = 50`nDate Created = July, 14, 2000`nminimum_balance = 100"
$ws.Range("G11").Value = "service_charge = 1"

# Test case 6 (row 12): __str__ / appropriate value returned based on attribute values.
$ws.Range("G12").Value = "Account Number 910 Balance: `$575.00`nMinimum Balance: `$100.00 Account Type: Savings"

# Rows of filled-in test cases grew to fit the multi-line inputs/outputs.
$ws.Rows.Item(7).RowHeight = 84
$ws.Rows.Item(8).RowHeight = 84
$ws.Rows.Item(9).RowHeight = 84
$ws.Rows.Item(10).RowHeight = 84
$ws.Rows.Item(11).RowHeight = 84
$ws.Rows.Item(12).RowHeight = 84

# View state left by the author when they finished editing.
$ws.Range("G12").Select()
$excel.ActiveWindow.Zoom = 80
